# AFDP-2327 Added 'Fulfill' to the list of possible next queues to all
# rules that have 'Fulfill' queue as a default return queue.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 25-31 in column E ("List of possible next queues") belong to rules
# whose "Default return queue" (column G) is "Fulfill". Append ",Fulfill"
# to the existing list of possible next queues for each of those rows.
25..31 | ForEach-Object {
    $row = $_
    $cell = $ws.Cells.Item($row, 5)   # column E
    $cell.Value2 = "$($cell.Value2),Fulfill"
}

$ws.Range("E31").Select() | Out-Null
